$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 13 (old orphan "5840671 - Francisco José Moreira Chaves" row with
# blank column A). This shifts rows 14-22 up to 13-21, dimension becomes A1:C21.
$ws.Rows.Item(13).Delete()

# After the shift, the text in column B/C no longer lines up with the text in
# column A the way it should, so every B/C cell from row 10 down needs to be
# re-pointed at its correct final value.

# Row 10 (Objetivos:) -> now shows the professor's name
$ws.Range("B10").Value = "5840671 - Francisco José Moreira Chaves"
$ws.Range("C10").Value = "5840671 - Francisco José Moreira Chaves"

# Row 13 (Programa resumido:) -> "Semestral"
$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"

# Row 15 (Programa:) -> "01/01/2018" (kept as literal text, not parsed as a
# date; leading apostrophe forces text, then we restore the original
# per-column formatting so the style index matches the unedited columns).
$ws.Range("B15").Value = "'01/01/2018"
$ws.Range("B15").Style = "Normal"
$ws.Range("B15").WrapText = $true
$ws.Range("B15").VerticalAlignment = -4160

$ws.Range("C15").Value = "'01/01/2018"
$ws.Range("C15").Style = "Normal"
$ws.Range("C15").WrapText = $true
$ws.Range("C15").VerticalAlignment = -4160
$ws.Range("C15").Font.Color = $ws.Range("C19").Font.Color

# Row 18 (Método:) -> "5840671 - Francisco José Moreira Chaves"
$ws.Range("B18").Value = "5840671 - Francisco José Moreira Chaves"
$ws.Range("C18").Value = "5840671 - Francisco José Moreira Chaves"

# Row 19 (Critério:) -> "Duas Provas  P1  1º bimestre e P2  2º bimestre"
$ws.Range("B19").Value = "Duas Provas  P1  1º bimestre e P2  2º bimestre"
$ws.Range("C19").Value = "Duas Provas  P1  1º bimestre e P2  2º bimestre"

# Row 20 (Norma de recuperação:) -> "MF = (P1+ P2)/2"
$ws.Range("B20").Value = "MF = (P1+ P2)/2"
$ws.Range("C20").Value = "MF = (P1+ P2)/2"

# Row 21 (Bibliografia:) -> "NF = (MF + PR)/ 2 , onde PR é uma prova de recuperação"
$ws.Range("B21").Value = "NF = (MF + PR)/ 2 , onde PR é uma prova de recuperação"
$ws.Range("C21").Value = "NF = (MF + PR)/ 2 , onde PR é uma prova de recuperação"

# Dimension will recompute automatically; worksheet now spans A1:C21
